$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.884.74'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.15%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.468.82'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.99%  '

$ws.Range("E4").Value = '  -0.29%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.97%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.37'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -5.23%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.466.89'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.93%  '

$ws.Range("E8").Value = '  -0.35%  '

$ws.Range("E9").Value = '  -1.33%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.143'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.44%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '7.56'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.84%  '

$ws.Range("E12").Value = '  -2.82%  '

$ws.Range("E13").Value = '  -3.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '31.84'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.59%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.050.36'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.459.26'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.51%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.894.32'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.05%  '

$ws.Range("E18").Value = '  -0.49%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.46'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.14%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.41'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.43%  '

$ws.Range("E21").Value = '  -0.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '439.81'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -4.30%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.612'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.21%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '78.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.63%  '

$ws.Range("E25").Value = '  +0.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.602.36'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.00%  '

$ws.Range("E27").Value = '  -8.71%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.82'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.42%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.47'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -8.07%  '

$ws.Range("E30").Value = '  -5.06%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.61'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.24%  '

$ws.Range("E32").Value = '  -2.07%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '25.44'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.14%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.09'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -6.00%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.456.35'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.90%  '

$ws.Range("E37").Value = '  -6.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.94'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -5.02%  '

$ws.Range("E39").Value = '  -0.04%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '173.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.63%  '

$ws.Range("E42").Value = '  -2.39%  '

$ws.Range("E43").Value = '  -9.32%  '

$ws.Range("E44").Value = '  -3.78%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.887'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.83%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.49'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -6.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '46.10'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.60%  '

$ws.Range("E48").Value = '  -8.57%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.49'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -8.86%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.48'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.60%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.991'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -4.02%  '
